$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: insert a new "Send Email" checkbox column before column A ---
$ws1.Columns.Item(1).Insert()
$ws1.Columns.Item(1).ColumnWidth = 13

# Headers (write "Email Address" before "Send Email" so new shared-string
# entries land in the same order as the target file).
$ws1.Range("B1").Value = "Email Address"
$ws1.Range("A1").Value = "Send Email"

# New boolean column values
$ws1.Range("A2").Value = $true
$ws1.Range("A3").Value = $false

# Hyperlinks got orphaned by the column insert (engine doesn't auto-shift
# them) - drop the stale ones and re-add against the new column B, in the
# same order as the target (B3 then B2).
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B3"), "mailto:vwm5@cdc.gov")
$ws1.Range("B3").Style = "Hyperlink"
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:lamtahri@hotmail.com")
$ws1.Range("B2").Style = "Hyperlink"

# AutoFilter across the full header row
$ws1.Range("A1:K1").AutoFilter()
$n = $ws1.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$K`$1")
$n.Visible = $false

# --- Sheet2: sender address update ---
$ws2.Range("C2").Value = "epiinfo@cdc.gov"
$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:epiinfo@cdc.gov")
$ws2.Range("C2").Style = "Hyperlink"
$ws2.Hyperlinks.Add($ws2.Range("C3"), "mailto:epiinfo@cdc.gov")
$ws2.Range("C3").Style = "Hyperlink"
$ws2.Range("C3").Value = "epiinfo@cdc.gov"

# --- Active sheet / selections swap from Sheet2 to Sheet1 ---
$ws1.Activate()
$ws1.Range("K25").Select()
$ws2.Range("C7").Select()
$ws1.Activate()

Write-Output "done"
